$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 154 - this shifts the existing rows 154-210
# down to 155-211 (format is carried down from the row above, matching
# the D column's date style).
$ws.Rows.Item(154).Insert()

# Populate the new row 154 with the new record. The categorical columns
# (A,B,C,E,F,G,H,I,J,K,L) repeat the same values as the rest of this
# "Femacal de La Calera" / Arandano (blue) block.
$ws.Cells.Item(154, 1).Value = 3
$ws.Cells.Item(154, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(154, 3).Value = "Coquimbo"
$ws.Cells.Item(154, 4).Value = 44636
$ws.Cells.Item(154, 5).Value = 5
$ws.Cells.Item(154, 6).Value = "Fruta"
$ws.Cells.Item(154, 7).Value = 100101
$ws.Cells.Item(154, 8).Value = "Berries"
$ws.Cells.Item(154, 9).Value = 100101001
$ws.Cells.Item(154, 10).Value = "Arándano (blue)"
$ws.Cells.Item(154, 11).Value = "Sin especificar"
$ws.Cells.Item(154, 12).Value = "Primera"
$ws.Cells.Item(154, 13).Value = 25
$ws.Cells.Item(154, 14).Value = 4000
$ws.Cells.Item(154, 15).Value = 4000
$ws.Cells.Item(154, 16).Value = 4000
$ws.Cells.Item(154, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(154, 18).Value = "Provincia de Linares"
$ws.Cells.Item(154, 19).Value = 2000
$ws.Cells.Item(154, 20).Value = 2
